$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (longTraining) description text and average age
$ws.Range("F9").Value = "Government sponsored Training Program lasting longer than 6 months"
$ws.Range("E9").Value = 35.2

# Row 10 - shortTraining
$ws.Range("A10").Value = "shortTraining"
$ws.Range("B10").Value = "Short Training"
$ws.Range("C10").Value = 1993
$ws.Range("D10").Value = "Active Labor Market Policy"
$ws.Range("E10").Value = 34.8
$ws.Range("F10").Value = "Government sponsored Training Program lasting less less then 6 months"

# Row 11 - practiceFirm
$ws.Range("A11").Value = "practiceFirm"
$ws.Range("B11").Value = "Practice Firm"
$ws.Range("C11").Value = 1993
$ws.Range("D11").Value = "Active Labor Market Policy"
$ws.Range("E11").Value = 36

# Row 12 - retraining (description filled before its own name/label, and before practiceFirm's F)
$ws.Range("C12").Value = 1993
$ws.Range("D12").Value = "Active Labor Market Policy"
$ws.Range("E12").Value = 31
$ws.Range("F12").Value = "Government sponsored training to obtain a new professional degree in a field other than the`nprofession currently held"

$ws.Range("F11").Value = "Government sponsored training that simulates a job"

$ws.Range("A12").Value = "retraining"
$ws.Range("B12").Value = "Retraining"

# Row 13 - classRoomTraining
$ws.Range("A13").Value = "classRoomTraining"
$ws.Range("B13").Value = "Class Room Training"
$ws.Range("C13").Value = 2000
$ws.Range("D13").Value = "Active Labor Market Policy"
$ws.Range("E13").Value = 37

# Row 14 - shortTrainingPost1998 / Trainingsmaßnahmen
$ws.Range("B14").Value = "Trainingsmaßnahmen"
$ws.Range("C14").Value = 2000
$ws.Range("D14").Value = "Active Labor Market Policy"
$ws.Range("E14").Value = 37
$ws.Range("F14").Value = "Government sponsored short-term training programs `"Trainingsmaßnahmen`", which where introduced again in 1998 after being abolished in 1993"

$ws.Range("A14").Value = "shortTrainingPost1998"

$ws.Range("F13").Value = "Government sponsored training that lasts on average 7.5 months."

# Apply wrap-text style (style index 1 in original) to F column cells that should wrap
$ws.Range("F10").WrapText = $true
$ws.Range("F12").WrapText = $true
$ws.Range("F13").WrapText = $true
$ws.Range("F14").WrapText = $true
$ws.Range("F11").Style = "Standard"

# View/selection changes
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("D14").Select()
